$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 11:40"

# Re-sort: swap adjacent country label pairs (Polonia/Portugal, Islas Feroe/Guadalupe, Islas Malvinas/Montserrat)
$ws.Range("A48").Value = "Polonia"
$ws.Range("A49").Value = "Portugal"
$ws.Range("A174").Value = "Islas Feroe"
$ws.Range("A175").Value = "Guadalupe"
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# Update statistics values
# Row 6
$ws.Range("B6").Value = 2271034
$ws.Range("C6").Value = 3881
$ws.Range("D6").Value = 1583652
$ws.Range("E6").Value = 641999
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 45383

# Row 18
$ws.Range("B18").Value = 263503
$ws.Range("C18").Value = 2996
$ws.Range("D18").Value = 151972
$ws.Range("E18").Value = 108060
$ws.Range("G18").Value = 33
$ws.Range("H18").Value = 3471

# Row 26
$ws.Range("B26").Value = 128776
$ws.Range("C26").Value = 1693
$ws.Range("D26").Value = 83710
$ws.Range("E26").Value = 39242
$ws.Range("G26").Value = 59
$ws.Range("H26").Value = 5824

# Row 33
$ws.Range("B33").Value = 85354
$ws.Range("C33").Value = 632
$ws.Range("D33").Value = 60019
$ws.Range("E33").Value = 24716
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 619

# Row 37
$ws.Range("B37").Value = 82050
$ws.Range("C37").Value = 263
$ws.Range("D37").Value = 76720
$ws.Range("E37").Value = 4797
$ws.Range("G37").Value = 12
$ws.Range("H37").Value = 533

# Row 48
$ws.Range("B48").Value = 52961
$ws.Range("C48").Value = 551
$ws.Range("D48").Value = 37150
$ws.Range("E48").Value = 13990
$ws.Range("G48").Value = 12
$ws.Range("H48").Value = 1821

# Row 49
$ws.Range("B49").Value = 52825
$ws.Range("D49").Value = 38600
$ws.Range("E49").Value = 12466
$ws.Range("H49").Value = 1759

# Row 57
$ws.Range("B57").Value = 37269
$ws.Range("C57").Value = 107
$ws.Range("D57").Value = 26415
$ws.Range("E57").Value = 9510
$ws.Range("G57").Value = 16
$ws.Range("H57").Value = 1344

# Row 71
$ws.Range("B71").Value = 22245
$ws.Range("C71").Value = 139
$ws.Range("D71").Value = 20123
$ws.Range("E71").Value = 1399

# Row 73
$ws.Range("B73").Value = 21269
$ws.Range("C73").Value = 397
$ws.Range("D73").Value = 9875
$ws.Range("E73").Value = 10824

# Row 79
$ws.Range("E79").Value = 6363
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 102

# Row 88
$ws.Range("B88").Value = 9103
$ws.Range("C88").Value = 9
$ws.Range("D88").Value = 8809
$ws.Range("E88").Value = 169

# Row 95
$ws.Range("B95").Value = 7623
$ws.Range("C95").Value = 22
$ws.Range("E95").Value = 310

# Row 111
$ws.Range("B111").Value = 4182
$ws.Range("C111").Value = 33
$ws.Range("D111").Value = 3052
$ws.Range("E111").Value = 1072

# Row 123
$ws.Range("B123").Value = 2615
$ws.Range("C123").Value = 16
$ws.Range("D123").Value = 1874
$ws.Range("E123").Value = 710

# Row 129
$ws.Range("B129").Value = 2272
$ws.Range("C129").Value = 17
$ws.Range("E129").Value = 183
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 129

# Row 174
$ws.Range("B174").Value = 318
$ws.Range("C174").Value = 12
$ws.Range("D174").Value = 225
$ws.Range("E174").Value = 93
$ws.Range("H174").Value = 0

# Row 175
$ws.Range("B175").Value = 317
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 186
$ws.Range("E175").Value = 117
$ws.Range("H175").Value = 14

# Row 213
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
